# fix bugs and add toefl mode
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new TOEFL-mode columns (Structure, Listening, Reading, Total)
$ws.Range("E1").Value = "Structure"
$ws.Range("F1").Value = "Listening"
$ws.Range("G1").Value = "Reading"
$ws.Range("H1").Value = "Total"

$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 300

# Move the active selection, matching the saved workbook state
[void]$ws.Range("C8").Select()
